# Feature selection based on threshold
# Remove rows whose "name" label is 2-46, 3-45, 2-52, 3-17, or 3-5
# (rows with missing / outlier data), leaving the rest of the table intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labelsToRemove = @("2-46", "3-45", "2-52", "3-17", "3-5")

# Find the last used row in column A
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

# Walk from the bottom up so row indices of rows above stay valid
for ($r = $lastRow; $r -ge 1; $r--) {
    $label = $ws.Cells.Item($r, 1).Value()
    if ($labelsToRemove -contains $label) {
        $ws.Rows.Item($r).Delete()
    }
}

# Match the resulting selection recorded in the saved workbook
$ws.Range("B14").Select()
